$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# New rows 28-33 on Sheet6. Cell values are written in the same order the
# original author entered them so that brand-new shared strings land at the
# same indices as in the target workbook.
$ws.Range("A28").Value = "Sales"
$ws.Range("B28").Value = 100

$ws.Range("A29").Value = "Step 1: Formula"
$ws.Range("D29").Value = "Step 2: N...C"

$ws.Range("A32").Value = "subnetmask?"
$ws.Range("A33").Value = "255.255.255.128"

$ws.Range("D30").Value = "nccc cccc"
$ws.Range("D31").Value = "0000 0000"
$ws.Range("D32").Value = "1000 0000"

$ws.Range("F32").Value = "192.168.0.128"
$ws.Range("F33").Value = "192.168.0.255"

# Make Sheet6 the active sheet/tab (was Sheet2) and restore the author's
# scroll position + selection on it.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A37").Select() | Out-Null
